$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G317").Value = 82
$ws.Range("G318").Value = 81
$ws.Range("G319").Value = 82
$ws.Range("G320").Value = 86
$ws.Range("G321").Value = 86
$ws.Range("G322").Value = 99
$ws.Range("G323").Value = 100
$ws.Range("G324").Value = 115
$ws.Range("G325").Value = 124
$ws.Range("G326").Value = 117

$ws.Range("C326").Value = 57
$ws.Range("C327").Value = 50
$ws.Range("C328").Value = 112

$ws.Range("G327").Value = 110
$ws.Range("G328").Value = 111

$ws.Range("L327").Value = 2
$ws.Range("M327").Value = 1

$ws.Range("C329").Value = 13
$ws.Range("E329").Value = 12
$ws.Range("F329").Value = 11
$ws.Range("G329").Value = 125
$ws.Range("L329").Value = 0
$ws.Range("M329").Value = 0
